# Daily attendance processing - 2025-12-01 07:31:28
# Reorders the "Recorded By" email lists for several sessions, updates a
# couple of attendance counts/percentages that changed after reprocessing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - ANATOMY session 1 recorders reordered
$ws.Range("G2").Value = "gehanadel@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# Row 3 - ANATOMY session 2 recorders reordered
$ws.Range("G3").Value = "majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, System, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# Row 4 - ANATOMY session 3 recorders reordered
$ws.Range("G4").Value = "majorelle.magdy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# Row 5 - ANATOMY session 4 recorders reordered + updated (added Amira.Sobhy) and attendance count updated
$ws.Range("G5").Value = "Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("H5").Value = "102/251"

# Row 6 - ANATOMY session 5 recorders reordered
$ws.Range("G6").Value = "majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg"

# Row 7 - BIOCHEMISTRY LAB/CBL session 1 recorders reordered
$ws.Range("G7").Value = "AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg"

# Row 9 - HISTOLOGY session 1 recorders reordered
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

# Row 10 - Average Attendance % stat updated (text, keep original formatting/style).
# Setting .Value directly would make Excel auto-detect the percent-looking
# text and silently convert the cell to a numeric percentage, changing both
# its stored type and its style, so the format is forced to Text first and
# the original cell formatting (fill/alignment) is restored afterwards from
# an identically-styled neighbour cell.
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "26.0%"
$ws.Range("K10").Copy() | Out-Null
$ws.Range("L10").PasteSpecial(-4122) | Out-Null

# Row 12 - MICROBIOLOGY session 1 recorders reordered
$ws.Range("G12").Value = "Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, dina.adel@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"

# Row 15 - PARASITOLOGY group statistics, Avg Attendance % updated (same trick as L10)
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "26.0%"
$ws.Range("R15").Copy() | Out-Null
$ws.Range("S15").PasteSpecial(-4122) | Out-Null

# Row 27 - PHARMACOLOGY session 2 recorders reordered
$ws.Range("G27").Value = "nourhan.mostafa@med.asu.edu.eg, hana.amr@med.asu.edu.eg"

# Row 28 - PHYSIOLOGY session 1 recorders reordered
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"

$excel.CutCopyMode = $false
